$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A72").Value = "global_get-free-chips-modal_header"
$ws.Range("A73").Value = "global_get-free-chips-modal_content"
$ws.Range("A74").Value = "global_get-free-chips-modal_btn-txt"
$ws.Range("B72").Value = "Refuel  "
$ws.Range("B74").Value = "Get Your Free Chips"
$ws.Range("C72").Value = "Auftanken"
$ws.Range("C74").Value = "Gratis Chips Holen"
$ws.Range("C73").Value = "Oh nein, es scheint, als würden Ihnen die Chips ausgehen! Aber keine Sorge, hier ist eine neue Charge Chips für Sie, damit Sie weiterspielen können!"
$ws.Range("B73").Value = "Oh noes, it seems like you're running out of chips! But don't worry, here's a fresh batch of chips for you so you can continue playing!"
